$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column (H) matching the formatting of the other
# header cells (copy format from G1, which uses the bold/border/centered
# header style) then set the new header text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding data value for row 2.
$ws.Range("H2").Value = 1
